$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update email addresses for Alice, Bob and Eve and turn them into mailto
# hyperlinks (matches the author's "sendmail" edit - real addresses were
# swapped for placeholder sina/163 addresses before mailing).
$ws.Range("B2").Value = "xxxxxxx@163.com"
$ws.Range("B3").Value = "XXXX@sina.com"
$ws.Range("B6").Value = "XXXXXXX@sina.com"

$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:xxxxxxx@163.com")
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:XXXXXXX@sina.com")
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:XXXX@sina.com")

# Bob's row (row 3) picked up a stray "paid" mark in column H.
$ws.Range("H3").Value = "paid"
